$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.724.79"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "3.053.56"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D5").Value = "'555.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "'144.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.14%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "3.047.27"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").Value = "'0.499"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'6.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.150"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "'0.466"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.90%  "
$ws.Range("D13").Value = "'0.0000226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").Value = "'34.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "3.572.25"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "63.717.68"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "3.059.46"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "'6.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").Value = "'473.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("D21").Value = "'13.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").Value = "'0.672"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.19%  "
$ws.Range("D23").Value = "'7.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.22%  "
$ws.Range("D24").Value = "'13.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.29%  "
$ws.Range("D25").Value = "'81.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'2.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'8.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "'2.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").Value = "'25.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").Value = "'1.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").Value = "'2.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").Value = "'5.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.94%  "
$ws.Range("D35").Value = "'6.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("D36").Value = "'54.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("D37").Value = "'458.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("D38").Value = "'2.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.53%  "
$ws.Range("D39").Value = "'0.0826"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").Value = "'0.0402"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.83%  "
$ws.Range("D41").Value = "2.943.78"
$ws.Range("E41").Value = "  -7.02%  "
$ws.Range("D42").Value = "'8.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").Value = "'0.113"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.73%  "
$ws.Range("D44").Value = "'27.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("D45").Value = "'0.257"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "'2.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("D48").Value = "'0.111"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").Value = "'119.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.28%  "
$ws.Range("D50").Value = "0.0₃0512"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("E51").Value = "  -1.40%  "
